# July07 Started Explore Section: Recursion #1
# P206, P344, P700
# P54 NF

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New date row for July 7, 2020 (one day after the existing 44018 / July 6 row)
$ws.Range("A3").Value = 44019

# Build the "Explore Section" tracker list in column F (rows 3-11),
# plus a couple of blank placeholder rows below it (12-14) and an
# empty leading cell at F2 so the block lines up under the header row.
# (Touching Font.Bold with a no-op keeps the blank cells materialized
# in the sheet, matching the formatted-but-empty placeholder cells.)
$ws.Range("F2").Font.Bold = $False
$ws.Range("F3").Value = "Explore Section (Sequence as Below)"
$ws.Range("F4").Value = "Recursion 1"
$ws.Range("F5").Value = "Recursion 2"
$ws.Range("F6").Value = "Binary Search"
$ws.Range("F7").Value = "Binary Tree"
$ws.Range("F8").Value = "Binary Search Tree"
$ws.Range("F9").Value = "HashTable"
$ws.Range("F10").Value = "Array and String"
$ws.Range("F11").Value = "Linked List"
$ws.Range("F12").Font.Bold = $False
$ws.Range("F13").Font.Bold = $False
$ws.Range("F14").Font.Bold = $False

# Note that we've started the Explore section in the main log.
$ws.Range("C3").Value = "Start Explore section"

# Leave the selection where the user's cursor ended up.
$ws.Range("C4").Select() | Out-Null
